$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert 3 new blank rows before row 38 ---
# This pushes the existing rows 38 and 39 (Black Amber / Primera, Black Amber /
# Segunda) down to rows 41 and 42 automatically, content intact.
$ws.Rows("38:40").Insert()

# --- 2) Fill the newly inserted row 38 with what used to be row 35's data
#        (Angeleno / Primera, fecha 44615, Región Metropolitana) ---
$ws.Range("A38").Value = 2
$ws.Range("B38").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44615
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = "Frutos de hueso (carozo)"
$ws.Range("I38").Value = 100103002
$ws.Range("J38").Value = "Ciruela"
$ws.Range("K38").Value = "Angeleno"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 16
$ws.Range("N38").Value = 200000
$ws.Range("O38").Value = 210000
$ws.Range("P38").Value = 205000
$ws.Range("Q38").Value = "`$/bins (450 kilos)"
$ws.Range("R38").Value = "Región Metropolitana"
$ws.Range("S38").Value = 456
$ws.Range("T38").Value = 450

# --- 3) Fill the newly inserted row 39 with what used to be row 36's data
#        (Angeleno / Segunda, fecha 44615, Región Metropolitana) ---
$ws.Range("A39").Value = 2
$ws.Range("B39").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44615
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103002
$ws.Range("J39").Value = "Ciruela"
$ws.Range("K39").Value = "Angeleno"
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 20
$ws.Range("N39").Value = 160000
$ws.Range("O39").Value = 170000
$ws.Range("P39").Value = 165000
$ws.Range("Q39").Value = "`$/bins (450 kilos)"
$ws.Range("R39").Value = "Región Metropolitana"
$ws.Range("S39").Value = 367
$ws.Range("T39").Value = 450

# --- 4) Fill the newly inserted row 40 with what used to be row 37's data
#        (Black Amber / Especial, fecha 44595, Región de O'Higgins) ---
$ws.Range("A40").Value = 2
$ws.Range("B40").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44595
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = "Frutos de hueso (carozo)"
$ws.Range("I40").Value = 100103002
$ws.Range("J40").Value = "Ciruela"
$ws.Range("K40").Value = "Black Amber"
$ws.Range("L40").Value = "Especial"
$ws.Range("M40").Value = 160
$ws.Range("N40").Value = 15500
$ws.Range("O40").Value = 16000
$ws.Range("P40").Value = 15750
$ws.Range("Q40").Value = "`$/caja 15 kilos granel"
$ws.Range("R40").Value = "Región de O'Higgins"
$ws.Range("S40").Value = 1050
$ws.Range("T40").Value = 15

# --- 5) Update row 35 in place: becomes Angeleno / Especial, fecha 44644,
#        Región de O'Higgins ---
$ws.Range("D35").Value = 44644
$ws.Range("L35").Value = "Especial"
$ws.Range("M35").Value = 10
$ws.Range("N35").Value = 230000
$ws.Range("O35").Value = 240000
$ws.Range("P35").Value = 235000
$ws.Range("R35").Value = "Región de O'Higgins"
$ws.Range("S35").Value = 522

# --- 6) Update row 36 in place: becomes Angeleno / Primera, fecha 44644,
#        Región de O'Higgins ---
$ws.Range("D36").Value = 44644
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 16
$ws.Range("N36").Value = 210000
$ws.Range("O36").Value = 220000
$ws.Range("P36").Value = 215000
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 478

# --- 7) Update row 37 in place: becomes Angeleno / Segunda, fecha 44644,
#        unidad $/bins (450 kilos), Región Metropolitana ---
$ws.Range("D37").Value = 44644
$ws.Range("K37").Value = "Angeleno"
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 16
$ws.Range("N37").Value = 150000
$ws.Range("O37").Value = 160000
$ws.Range("P37").Value = 155000
$ws.Range("Q37").Value = "`$/bins (450 kilos)"
$ws.Range("S37").Value = 344
$ws.Range("T37").Value = 450
